$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds text-formatted numbers (e.g. "1.002", "22.035.87").
# Force text format first so Excel doesn't auto-coerce these into numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 16 <-> Row 17 swap: Chainlink and WrappedEther swap ranking positions
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.556.23"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "6.594"
$ws.Range("E17").Value = "  -2.41%  "

# Price / Volume(1h) updates for remaining rows
$ws.Range("D2").Value = "22.035.87"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.552.54"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "289.82"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "0.3928"
$ws.Range("E7").Value = "  +3.54%  "
$ws.Range("D8").Value = "0.3204"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "43.53"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "0.07218"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  -6.14%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "5.639"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").Value = "18.60"
$ws.Range("E14").Value = "  -7.89%  "
$ws.Range("D15").Value = "0.00001123"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D18").Value = "0.06586"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "83.14"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "6.264"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "15.39"
$ws.Range("E22").Value = "  -4.96%  "
$ws.Range("D23").Value = "11.21"
$ws.Range("E23").Value = "  -4.32%  "
$ws.Range("D24").Value = "22.048.58"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "2.377"
$ws.Range("E25").Value = "  +3.77%  "
$ws.Range("D26").Value = "2.410"
$ws.Range("E26").Value = "  -6.11%  "
$ws.Range("D27").Value = "148.26"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").Value = "18.45"
$ws.Range("E28").Value = "  -4.45%  "
$ws.Range("D29").Value = "4.873"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").Value = "1.726.05"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "118.14"
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").Value = "0.9884"
$ws.Range("E32").Value = "  -8.46%  "
$ws.Range("D33").Value = "5.750"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").Value = "0.08274"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").Value = "1.609"
$ws.Range("E35").Value = "  -16.30%  "
$ws.Range("D36").Value = "8.997"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("D37").Value = "0.02250"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "0.06032"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("D39").Value = "5.078"
$ws.Range("E39").Value = "  -5.39%  "
$ws.Range("D40").Value = "1.207"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "0.2033"
$ws.Range("E41").Value = "  -5.84%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "10.59"
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("D44").Value = "0.5776"
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("D45").Value = "3.738"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "12.92"
$ws.Range("E46").Value = "  -5.85%  "
$ws.Range("D47").Value = "0.5551"
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("D48").Value = "117.48"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "1.883"
$ws.Range("D50").Value = "1.127"
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("D51").Value = "0.06816"
$ws.Range("E51").Value = "  -3.60%  "
